# Swap the data (columns B:AD) between pairs of adjacent rows.
# Column A (the row index) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(25, 26),
    @(39, 40),
    @(58, 59),
    @(74, 75),
    @(82, 83),
    @(84, 85),
    @(122, 123),
    @(143, 144),
    @(148, 149),
    @(151, 152),
    @(160, 161),
    @(171, 172)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $values1 = $range1.Value()
    $values2 = $range2.Value()

    $range1.Value = $values2
    $range2.Value = $values1
}
